$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = 131307610
$ws.Range("B28").Value = 99015
$ws.Range("D28").Value = "VU"
$ws.Range("E28").Value = 220787
$ws.Range("F28").Value = "Knärot"
$ws.Range("G28").Value = "Goodyera repens"
$ws.Range("H28").Value = "(L.) R. Br."

$ws.Range("I28").NumberFormat = "@"
$ws.Range("I28").Value = "130"
$ws.Range("I28").Style = "Normal"

$ws.Range("J28").Value = "stjälkar/strån/skott"
$ws.Range("K28").Value = "fullt utvecklade blad"
$ws.Range("P28").Value = "Bäckmossen, Vstm"
$ws.Range("Q28").Value = 521051
$ws.Range("R28").Value = 6625617
$ws.Range("S28").Value = 78
$ws.Range("T28").Value = "Örebro"
$ws.Range("U28").Value = "Lindesberg"
$ws.Range("V28").Value = "Västmanland"
$ws.Range("W28").Value = "Ramsberg"
$ws.Range("X28").Value = "T-Lin-0017"

$ws.Range("Y28").NumberFormat = "@"
$ws.Range("Y28").Value = "2025-12-02"
$ws.Range("Y28").Style = "Normal"

$ws.Range("AA28").NumberFormat = "@"
$ws.Range("AA28").Value = "2025-12-02"
$ws.Range("AA28").Style = "Normal"

$ws.Range("AD28").Value = $false
$ws.Range("AE28").Value = $false
$ws.Range("AG28").Value = $false

$ws.Range("AT28").NumberFormat = "@"
$ws.Range("AT28").Value = " "
$ws.Range("AT28").Style = "Normal"

$ws.Range("AW28").Value = "Sofia Lund"
$ws.Range("AX28").Value = "Lotta Sörman"
$ws.Range("AY28").Value = "Floraväkteri Sverige"
